$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A (everything currently in A/B/C shifts to B/C/D) ---
$ws.Columns("A").Insert()

# Carry the header row's highlight formatting (fillId=3, applied to B3:D3) onto
# the new A3 cell before we put the "S/N" label in it.
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New "S/N" header + row numbers in the freshly inserted column A.
$ws.Range("A3").Value = "S/N"
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3

# Columns("A").Insert() shifts the cell grid but leaves the worksheet-level
# Hyperlinks collection pointing at the old (pre-shift) addresses, so the 3
# mailto hyperlinks need to be re-created at their new column-D locations.
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:ding@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:beh@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:charles@yahoo.com")

# Selection moved to F13 in the saved file.
$ws.Range("F13").Select() | Out-Null
